$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 17 (the totals row). This shifts the old
# row 17 (totals) down to row 18 and the old row 18 (footer) down to row 19.
$ws.Rows("17:17").Insert()

# Copy the formatting of the previous item row (row 16) onto the new row 17
# so the new item row 11 ("محلول ملح") looks like the other item rows.
$ws.Range("A16:Q16").Copy()
$ws.Range("A17:Q17").PasteSpecial(-4122)

# Recreate the merged cells for the new item row, matching the pattern used
# by every other item row (A:B, C:G, H:K, L:M, N:O).
$ws.Range("A17:B17").Merge()
$ws.Range("C17:G17").Merge()
$ws.Range("H17:K17").Merge()
$ws.Range("L17:M17").Merge()
$ws.Range("N17:O17").Merge()

# Populate the new item row (#11 - محلول ملح / Saline solution).
$ws.Range("A17").Value = 11
$ws.Range("C17").Value = "محلول ملح"
$ws.Range("H17").Value = "16:0"
$ws.Range("L17").Value = "0"
$ws.Range("N17").Value = "24.00"
$ws.Range("P17").Value = "24.0000"
$ws.Range("Q17").Value = "1:0"

# Update the grand-total cell (old row 17, now shifted to row 18).
$ws.Range("P18").Value = 344

# Update the generated-on timestamp in the footer (old row 18, now row 19).
$ws.Range("A19").Value = "Thursday, 14 August, 2025 10:26 AM"
